$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "7 complete" -> "11 complete" comment for FOCUS SAX 50 scale row (I25)
$ws.Range("I25").Value = "11 complete"

# 2. Rows 20-26: touch column B (Date Completed) leaving it blank but formatted with the
#    default/general style (matches the reference edit's explicit blank B cells).
$ws.Range("B20:B26").HorizontalAlignment = 1

# 3. Row 27: EVHP Holder Goose Neck to TEE Clip Adapter now has a completed date
$ws.Range("B27").HorizontalAlignment = 1
$ws.Range("B27").Value = "27-05-2018"

# 4. Row 28: EVHP Holder Goose Neck to Cannula Mount Adapter completed date changes from 26-05-2018 to 27-05-2018
$ws.Range("B28").Value = "27-05-2018"

# 5. Selection / view state changes
$ws.Range("A1").Select()
$ws.Range("B24").Select()
